$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# The task "CU Generar recibo de pago." (row 6) is finished:
# - Status (F6) moves from "En proceso" to "Hecho"
# - 2 hours were consumed on day 2 (N6)
$ws.Range("F6").Value = "Hecho"
$ws.Range("N6").Value = 2

# Leave the final selection on N6 to match the saved view state
$ws.Range("N6").Select()

$wb.Save()
